# Quiz result sheet: the marking engine previously choked on non-integer
# ("float") scoring input and silently left the student's answers/summary
# blank (marked "Absent"). This records the actual per-question student
# answers and fixes up the summary block now that scoring works again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# ---- Summary block (rows 9-12) -------------------------------------------
# Row/column header cells (A10/A11/A12) pick up the same boxed "mtitle"
# look already used by the A9 header cell above them.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

$ws.Range("B10").Value = 18    # Right
$ws.Range("C10").Value = 4     # Wrong
$ws.Range("D10").Value = 6     # Not Attempt
$ws.Range("E10").Value = 28    # Max

$ws.Range("B11").Value = 4     # Marking per right answer
$ws.Range("C11").Value = -1    # Marking per wrong answer (now numeric, not text)

$ws.Range("B12").Value = 72    # Total right marks (18 * 4)
$ws.Range("C12").Value = -4    # Total wrong marks (4 * -1)
$ws.Range("E12").Value = "68/112"

# ---- Per-question answers (rows 16-40) ------------------------------------
# The sheet has three side-by-side "Student Ans / Correct Ans" blocks
# (A:B, D:E, G:H). Only the first exam attempt (A:B) is kept; the other two
# stale blocks are removed entirely. The "Student Ans" column (A) is filled
# in with what the student actually answered, styled green ("correctStyle")
# when it matches the "Correct Ans" column (B) and red ("incorrectStyle")
# when it doesn't; rows left blank mean the question wasn't attempted.

function Set-Answer($cellRef, $answer, $isCorrect) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $answer
    if ($isCorrect) {
        $cell.Style = "correctStyle"
    } else {
        $cell.Style = "incorrectStyle"
    }
}

Set-Answer "A16" "Option A" $true
Set-Answer "A17" "Option D" $true
Set-Answer "A18" "Option B" $true
Set-Answer "A19" "Option C" $true
Set-Answer "A21" "Option C" $true
Set-Answer "A22" "Option D" $true
Set-Answer "A23" "Option A" $false
Set-Answer "A25" "Option A" $true
Set-Answer "A26" "Option B" $false
Set-Answer "A27" "Option A" $true
Set-Answer "A28" "Option D" $true
Set-Answer "A30" "Option B" $true
Set-Answer "A31" "Option D" $true
Set-Answer "A32" "Option C" $true
Set-Answer "A34" "Option C" $false
Set-Answer "A36" "Option A" $true
Set-Answer "A38" "Option C" $false
Set-Answer "A39" "Option D" $true
Set-Answer "A40" "Option D" $true
# Rows 20, 24, 29, 33, 35, 37 are left as-is (not attempted -> still blank).

# The second exam block (D:E) keeps its header (row 15) and first three
# rows, now filled in the same way as column A, but the remainder (rows
# 19-40) along with the whole third block (G:H, rows 15-21) is dropped.
Set-Answer "D16" "Option A" $true
Set-Answer "D17" "Option C" $true
Set-Answer "D18" "Option D" $true

$ws.Range("D19:E40").Clear()
$ws.Range("G15:H21").Clear()
